# Generate Report for Handback
# Update the generated timestamp values on the handback status report.
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 5875b730-... row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-16 02:41:42"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the 5875b730-... row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-16 02:41:37"
$wsZhCn.Range("K3").Value = "2016-08-16 02:42:07"

# de-de sheet: "Correspond Handback DateTime" for the
# 5875b730-... row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-08-16 02:42:15"
